$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 887.75
$ws.Range("I12").Value = 300.42856
$ws.Range("J12").Value = 4999
$ws.Range("K12").Value = 300.42856
$ws.Range("L12").Value = 4999
$ws.Range("M12").Value = -130.42856
$ws.Range("N12").Value = -5339

$ws.Range("H17").Value = 2703674.2
$ws.Range("J17").Value = 2941969
$ws.Range("L17").Value = 8825907
$ws.Range("N17").Value = -8826243

$ws.Range("H19").Value = 3389.8462
$ws.Range("I19").Value = 2761.4
$ws.Range("K19").Value = 2761.4
$ws.Range("M19").Value = -2586.4
$ws.Range("N19").ClearContents()

$ws.Range("H33").Value = 14702
$ws.Range("I33").Value = 16097.2
$ws.Range("K33").Value = 16097.2
$ws.Range("M33").Value = -15868.2
$ws.Range("N33").ClearContents()

$ws.Range("H40").Value = 5527
$ws.Range("I40").Value = 3557.8572
$ws.Range("J40").Value = 7250
$ws.Range("K40").Value = 3557.8572
$ws.Range("L40").Value = 7250
$ws.Range("M40").Value = -3382.8572
$ws.Range("N40").Value = -7600

$ws.Range("H80").Value = 276.2143
$ws.Range("I80").Value = 251.3077
$ws.Range("J80").Value = 600
$ws.Range("K80").Value = 753.9231
$ws.Range("L80").Value = 1800
$ws.Range("M80").Value = 244.0769
$ws.Range("N80").Value = -3796

$ws.Range("H83").Value = 276.2143
$ws.Range("I83").Value = 251.3077
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 2261.7693
$ws.Range("L83").Value = 5400
$ws.Range("M83").Value = 2730.2307
$ws.Range("N83").Value = -15384

$ws.Range("H103").Value = 201.8
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 201.8
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 605.4000000000001
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -1777.4

$ws.Range("H107").Value = 1044.6428
$ws.Range("I107").Value = 1111.5385
$ws.Range("J107").Value = 175
$ws.Range("K107").Value = 1111.5385
$ws.Range("L107").Value = 175
$ws.Range("M107").Value = 808.4614999999999
$ws.Range("N107").Value = -4015

$ws.Range("H135").Value = 440.1905
$ws.Range("I135").Value = 440.1905
$ws.Range("K135").Value = 3961.7145
$ws.Range("M135").Value = -1426.7145

$ws.Range("H138").Value = 2206.375
$ws.Range("I138").Value = 2031.25
$ws.Range("J138").Value = 2235.5625
$ws.Range("K138").Value = 6093.75
$ws.Range("L138").Value = 6706.6875
$ws.Range("M138").Value = -953.75
$ws.Range("N138").Value = -16986.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7576.4595
$ws.Range("I32").Value = 7648.794
$ws.Range("K32").Value = 7648.794
$ws.Range("M32").Value = -7361.794
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 256875
$ws.Range("J45").Value = 9250
$ws.Range("L45").Value = 9250
$ws.Range("N45").Value = -10004

$ws.Range("H61").Value = 7334.3
$ws.Range("I61").Value = 5981.5293
$ws.Range("K61").Value = 5981.5293
$ws.Range("M61").Value = -5769.5293
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 3769.476
$ws.Range("I74").Value = 3207.95
$ws.Range("K74").Value = 3207.95
$ws.Range("M74").Value = -2333.95
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 3769.476
$ws.Range("I77").Value = 3207.95
$ws.Range("K77").Value = 16039.75
$ws.Range("M77").Value = -11671.75
$ws.Range("N77").ClearContents()

$ws.Range("H101").Value = 21000
$ws.Range("I101").Value = 21000
$ws.Range("K101").Value = 21000
$ws.Range("M101").Value = -17755

$ws.Range("H132").Value = 4094
$ws.Range("I132").Value = 3228.647
$ws.Range("J132").Value = 11449.5
$ws.Range("K132").Value = 9685.940999999999
$ws.Range("L132").Value = 34348.5
$ws.Range("M132").Value = -7155.940999999999
$ws.Range("N132").Value = -39408.5

$ws.Range("H136").Value = 7334.3
$ws.Range("I136").Value = 5981.5293
$ws.Range("K136").Value = 17944.5879
$ws.Range("M136").Value = -15394.5879
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 170298.33
$ws.Range("J139").Value = 170298.33
$ws.Range("L139").Value = 170298.33
$ws.Range("N139").Value = -180578.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3731.75
$ws.Range("I94").Value = 3040.0625
$ws.Range("K94").Value = 3040.0625
$ws.Range("M94").Value = -2589.0625
$ws.Range("N94").ClearContents()

$ws.Range("H100").Value = 18083.75
$ws.Range("J100").Value = 18083.75
$ws.Range("L100").Value = 18083.75
$ws.Range("N100").Value = -20247.75

$ws.Range("H107").Value = 4991.3335
$ws.Range("I107").Value = 4691.533
$ws.Range("K107").Value = 4691.533
$ws.Range("M107").Value = -2771.533
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 5115.1665
$ws.Range("I134").Value = 4902.7827
$ws.Range("K134").Value = 14708.3481
$ws.Range("M134").Value = -12173.3481
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 455.07144
$ws.Range("I19").Value = 451.6154
$ws.Range("K19").Value = 451.6154
$ws.Range("M19").Value = -281.6154
$ws.Range("N19").ClearContents()

$ws.Range("H24").Value = 455.07144
$ws.Range("I24").Value = 451.6154
$ws.Range("K24").Value = 451.6154
$ws.Range("M24").Value = -281.6154
$ws.Range("N24").ClearContents()

$ws.Range("H31").Value = 9277.444
$ws.Range("I31").Value = 7249.6665
$ws.Range("J31").Value = 13333
$ws.Range("K31").Value = 7249.6665
$ws.Range("L31").Value = 13333
$ws.Range("M31").Value = -6954.6665
$ws.Range("N31").Value = -13923

$ws.Range("H34").Value = 9277.444
$ws.Range("I34").Value = 7249.6665
$ws.Range("J34").Value = 13333
$ws.Range("K34").Value = 7249.6665
$ws.Range("L34").Value = 13333
$ws.Range("M34").Value = -7047.6665
$ws.Range("N34").Value = -13737

$ws.Range("H58").Value = 12186.375
$ws.Range("I58").Value = 7497.3335
$ws.Range("J58").Value = 14999.8
$ws.Range("K58").Value = 7497.3335
$ws.Range("L58").Value = 14999.8
$ws.Range("M58").Value = -7294.3335
$ws.Range("N58").Value = -15405.8

$ws.Range("H109").Value = 49643
$ws.Range("J109").Value = 49643
$ws.Range("L109").Value = 49643
$ws.Range("N109").Value = -51723

$ws.Range("H136").Value = 12186.375
$ws.Range("I136").Value = 7497.3335
$ws.Range("J136").Value = 14999.8
$ws.Range("K136").Value = 22492.0005
$ws.Range("L136").Value = 44999.39999999999
$ws.Range("M136").Value = -19942.0005
$ws.Range("N136").Value = -50099.39999999999

$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13896513
$ws.Range("I4").Value = 12433561
$ws.Range("K4").Value = 37300683
$ws.Range("M4").Value = -37300571
$ws.Range("N4").ClearContents()

$ws.Range("H15").Value = 126096.836
$ws.Range("I15").Value = 1622.75
$ws.Range("J15").Value = 375045
$ws.Range("K15").Value = 4868.25
$ws.Range("L15").Value = 1125135
$ws.Range("M15").Value = -4728.25
$ws.Range("N15").Value = -1125415

$ws.Range("H86").Value = 414.2857
$ws.Range("I86").Value = 260
$ws.Range("J86").Value = 800
$ws.Range("K86").Value = 780
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = 406
$ws.Range("N86").Value = -4772

$ws.Range("H89").Value = 414.2857
$ws.Range("I89").Value = 260
$ws.Range("J89").Value = 800
$ws.Range("K89").Value = 2340
$ws.Range("L89").Value = 7200
$ws.Range("M89").Value = 3588
$ws.Range("N89").Value = -19056

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 61750
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 999
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 999
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 999
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1279

$ws.Range("H55").Value = 1036.4706
$ws.Range("I55").Value = 1072.6364
$ws.Range("J55").Value = 970.1667
$ws.Range("K55").Value = 1072.6364
$ws.Range("L55").Value = 970.1667
$ws.Range("M55").Value = -899.6364000000001
$ws.Range("N55").Value = -1316.1667

$ws.Range("H122").Value = 3034.1428
$ws.Range("J122").Value = 3049.3333
$ws.Range("L122").Value = 9147.999899999999
$ws.Range("N122").Value = -14047.9999

$ws.Range("H128").Value = 97993
$ws.Range("J128").Value = 97993
$ws.Range("L128").Value = 97993
$ws.Range("N128").Value = -107953

$ws.Range("H132").Value = 20849.455
$ws.Range("I132").Value = 25106.375
$ws.Range("K132").Value = 75319.125
$ws.Range("M132").Value = -72789.125
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 222222600
$ws.Range("I6").Value = 222222600
$ws.Range("K6").Value = 222222600
$ws.Range("M6").Value = -222222485

$ws.Range("H12").Value = 14000
$ws.Range("I12").Value = 14000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -13858
$ws.Range("N12").ClearContents()

$ws.Range("H132").Value = 5286.615
$ws.Range("I132").Value = 4893.8335
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 14681.5005
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -12151.5005
$ws.Range("N132").Value = -35060

$ws.Range("H136").Value = 3133.0312
$ws.Range("I136").Value = 2152.0356
$ws.Range("K136").Value = 6456.1068
$ws.Range("M136").Value = -3906.1068
$ws.Range("N136").ClearContents()
